# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold font, border, centered) from
# the last header cell (H1) onto the two new header cells so they match
# the look of the other headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-9: I column then J column values.
$values = @{
    2 = @(1, 3)
    3 = @(1, 6)
    4 = @(1, 6)
    5 = @(1, 5)
    6 = @(1, 4)
    7 = @(6, 8)
    8 = @(4, 6)
    9 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
